# Auto-generated-ish edit script: update Leve profit calc cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 355.83
$ws.Range("I15").Value = 355.83
$ws.Range("K15").Value = 1067.49
$ws.Range("M15").Value = -898.49

$ws.Range("H38").Value = 1000648
$ws.Range("I38").Value = 80
$ws.Range("J38").Value = 2501500
$ws.Range("K38").Value = 240
$ws.Range("L38").Value = 7504500
$ws.Range("M38").Value = 132
$ws.Range("N38").Value = -7505244

$ws.Range("H46").Value = 200
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -481
$ws.Range("N46").Value = $null

$ws.Range("H60").Value = 200
$ws.Range("I60").Value = 200
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -116
$ws.Range("N60").Value = $null

$ws.Range("H113").Value = 2670.4443
$ws.Range("I113").Value = 2243.8462
$ws.Range("J113").Value = 3779.6
$ws.Range("K113").Value = 2243.8462
$ws.Range("L113").Value = 3779.6
$ws.Range("M113").Value = 1010.1538
$ws.Range("N113").Value = -10287.6

$ws.Range("H128").Value = 60446.668
$ws.Range("J128").Value = 60446.668
$ws.Range("L128").Value = 60446.668
$ws.Range("N128").Value = -70406.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2391.8867
$ws.Range("I61").Value = 2311.422
$ws.Range("J61").Value = 2844.5
$ws.Range("K61").Value = 2311.422
$ws.Range("L61").Value = 2844.5
$ws.Range("M61").Value = -2099.422
$ws.Range("N61").Value = -3268.5

$ws.Range("H74").Value = 2041.6842
$ws.Range("I74").Value = 1819.2727
$ws.Range("J74").Value = 2347.5
$ws.Range("K74").Value = 1819.2727
$ws.Range("L74").Value = 2347.5
$ws.Range("M74").Value = -945.2727
$ws.Range("N74").Value = -4095.5

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

$ws.Range("H77").Value = 2041.6842
$ws.Range("I77").Value = 1819.2727
$ws.Range("J77").Value = 2347.5
$ws.Range("K77").Value = 9096.363499999999
$ws.Range("L77").Value = 11737.5
$ws.Range("M77").Value = -4728.363499999999
$ws.Range("N77").Value = -20473.5

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

$ws.Range("H92").Value = 40016.668
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40016.668
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40016.668
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -45008.668

$ws.Range("H122").Value = 6121.1787
$ws.Range("I122").Value = 6515.72
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 19547.16
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -17097.16
$ws.Range("N122").Value = -13399.9999

$ws.Range("H136").Value = 2391.8867
$ws.Range("I136").Value = 2311.422
$ws.Range("J136").Value = 2844.5
$ws.Range("K136").Value = 6934.266
$ws.Range("L136").Value = 8533.5
$ws.Range("M136").Value = -4384.266
$ws.Range("N136").Value = -13633.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

$ws.Range("H134").Value = 3617.102
$ws.Range("I134").Value = 2173
$ws.Range("J134").Value = 5121.375
$ws.Range("K134").Value = 6519
$ws.Range("L134").Value = 15364.125
$ws.Range("M134").Value = -3984
$ws.Range("N134").Value = -20434.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3130.6177
$ws.Range("I86").Value = 2772.5293
$ws.Range("J86").Value = 3488.7058
$ws.Range("K86").Value = 2772.5293
$ws.Range("L86").Value = 3488.7058
$ws.Range("M86").Value = -1649.5293
$ws.Range("N86").Value = -5734.7058

$ws.Range("H89").Value = 3130.6177
$ws.Range("I89").Value = 2772.5293
$ws.Range("J89").Value = 3488.7058
$ws.Range("K89").Value = 13862.6465
$ws.Range("L89").Value = 17443.529
$ws.Range("M89").Value = -8246.646500000001
$ws.Range("N89").Value = -28675.529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 771692.0600000001
$ws.Range("I2").Value = 108.181816
$ws.Range("J2").Value = 1984181
$ws.Range("K2").Value = 649.0908959999999
$ws.Range("L2").Value = 11905086
$ws.Range("M2").Value = -536.0908959999999
$ws.Range("N2").Value = -11905312

$ws.Range("H118").Value = 2453.625
$ws.Range("I118").Value = 757.25
$ws.Range("K118").Value = 2271.75
$ws.Range("M118").Value = -1028.75

$ws.Range("H122").Value = 1248.2084
$ws.Range("I122").Value = 1170.2667
$ws.Range("J122").Value = 1378.1111
$ws.Range("K122").Value = 10532.4003
$ws.Range("L122").Value = 12402.9999
$ws.Range("M122").Value = -8082.400299999999
$ws.Range("N122").Value = -17302.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 15760
$ws.Range("J76").Value = 15760
$ws.Range("L76").Value = 15760
$ws.Range("N76").Value = -16436

$ws.Range("H79").Value = 15760
$ws.Range("J79").Value = 15760
$ws.Range("L79").Value = 15760
$ws.Range("N79").Value = -18100

$ws.Range("H122").Value = 7276.357
$ws.Range("I122").Value = 8722.5
$ws.Range("K122").Value = 26167.5
$ws.Range("M122").Value = -23717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 37795
$ws.Range("J57").Value = 37795
$ws.Range("L57").Value = 37795
$ws.Range("N57").Value = -39303

$ws.Range("H75").Value = 38933.332
$ws.Range("J75").Value = 38933.332
$ws.Range("L75").Value = 38933.332
$ws.Range("N75").Value = -40805.332

$ws.Range("H78").Value = 38933.332
$ws.Range("J78").Value = 38933.332
$ws.Range("L78").Value = 116799.996
$ws.Range("N78").Value = -126159.996

$ws.Range("H132").Value = 1577.898
$ws.Range("I132").Value = 1415.5641
$ws.Range("J132").Value = 2211
$ws.Range("K132").Value = 4246.692300000001
$ws.Range("L132").Value = 6633
$ws.Range("M132").Value = -1716.692300000001
$ws.Range("N132").Value = -11693
